$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.38601601592701229
$ws.Range("B1").Value = 0.38516358302067033
$ws.Range("A2").Value = -0.33905641723548641
$ws.Range("B2").Value = 0.33549021101385712
$ws.Range("A3").Value = -0.15623979884758299
$ws.Range("B3").Value = 0.15566588763204336
$ws.Range("A4").Value = -0.1436658878360717
$ws.Range("B4").Value = 0.14316049269256581
$ws.Range("A5").Value = -0.13716049348805104
$ws.Range("B5").Value = 0.13615861009196095
$ws.Range("A6").Value = -0.0088045424901079805
$ws.Range("B6").Value = 0.0088003563261347395
$ws.Range("A7").Value = 0.011199642719351033
$ws.Range("B7").Value = -0.011205806884207803
$ws.Range("A8").Value = 0.031205805931305619
$ws.Range("B8").Value = -0.0312724387593315
$ws.Range("A9").Value = 0.037272437951453519
$ws.Range("B9").Value = -0.037369502916113184
$ws.Range("A10").Value = -0.02038577455457613
$ws.Range("B10").Value = 0.020386322638948684
$ws.Range("A11").Value = -0.051462841594897668
$ws.Range("B11").Value = 0.051387501675272063
$ws.Range("A12").Value = -0.045387502478767772
$ws.Range("B12").Value = 0.04514730672298306
$ws.Range("A13").Value = -0.03914730753836615
$ws.Range("B13").Value = 0.039082209872154827
$ws.Range("A14").Value = -0.027082210750532631
$ws.Range("B14").Value = 0.027051168773298606
$ws.Range("A15").Value = -0.021051169594633379
$ws.Range("B15").Value = 0.021026786117801066
$ws.Range("A16").Value = -0.015026786941523262
$ws.Range("B16").Value = 0.015003931110623814
$ws.Range("A17").Value = -0.0090039319375243565
$ws.Range("B17").Value = 0.0089999991426070736
$ws.Range("A18").Value = -0.03610716657248858
$ws.Range("B18").Value = 0.036096219258194395
$ws.Range("A19").Value = -0.027096220035255136
$ws.Range("B19").Value = 0.027013064763723449
$ws.Range("A20").Value = -0.018013065547908269
$ws.Range("B20").Value = 0.018004198611780708
$ws.Range("A21").Value = -0.0090041993970251255
$ws.Range("B21").Value = 0.0089999992140503693
$ws.Range("A22").Value = -0.093949889073487114
$ws.Range("B22").Value = 0.093636609109594815
$ws.Range("A23").Value = -0.084636609924173989
$ws.Range("B23").Value = 0.084127053080335301
$ws.Range("A24").Value = -0.04212705423589469
$ws.Range("B24").Value = 0.041999998837981778
$ws.Range("A25").Value = -0.12135406841737861
$ws.Range("B25").Value = 0.12116027981682009
$ws.Range("A26").Value = -0.061680537591985996
$ws.Range("B26").Value = 0.061549134046693155
$ws.Range("A27").Value = -0.05554913486858748
$ws.Range("B27").Value = 0.055111678623161797
$ws.Range("A28").Value = -0.049111679455989155
$ws.Range("B28").Value = 0.048825670566943202
$ws.Range("A29").Value = -0.036825671465955168
$ws.Range("B29").Value = 0.036706859731472363
$ws.Range("A30").Value = -0.016706860713123106
$ws.Range("B30").Value = 0.016676067189795862
$ws.Range("A31").Value = -0.0016760681245298059
$ws.Range("B31").Value = 0.0016687376927162489
$ws.Range("A32").Value = -0.0060005325676621979
$ws.Range("B32").Value = 0.0059999991544712472

$ws.Range("B1").ColumnWidth = 13.85
